# Apply the table style change recorded in the commit.
#
# The deck has a single table (the "B1 - TYPES OF FINANCIAL DOCUMENTS"
# slide, slide 5, shape 2). Its <a:tblPr> currently references the
# built-in table style {272ABCB1-C9E1-4D77-AD4F-41B4C0B512EE}. The author
# picked a different built-in table style from PowerPoint's Table Styles
# gallery, which re-stamps the table's <a:tableStyleId> with the new
# style's GUID: {CC8E522D-31DB-425B-AF21-0D2F99176C67}.
#
# PowerPoint's object model exposes this as Table.ApplyStyle(styleId) —
# Table.Style is read-only (attempting `Table.Style = ...` raises "Table
# styles cannot be assigned through a property - call Table.ApplyStyle
# instead"), so ApplyStyle is the correct call.

$p = $ppt.ActivePresentation

$oldStyleId = "{272ABCB1-C9E1-4D77-AD4F-41B4C0B512EE}"
$newStyleId = "{CC8E522D-31DB-425B-AF21-0D2F99176C67}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
